# "updated for 12 tarikh" — enter the Oct 12 (2021-10-12, serial 44481) row:
# bazar (market) charge for the day, and each person's lunch+dinner meal
# count. Every other changed cell in the sheet (U15, B18, F35, L35:U35,
# K36:T36, K37:T37, B38) is a formula that recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bazar charge entered for the day
$ws.Range("F15").Value = 1400

# Meal counts per person for the day (Antor, Ovi, Rajon, Shakib, Dhrubo,
# Gopal, Rahat, Roni, Shovon, Tawhid)
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 2
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = 2
$ws.Range("S15").Value = 0
$ws.Range("T15").Value = 2

# Move the active selection/cursor to where the user finished editing
[void]$ws.Range("D35").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
